$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Plot power output (single HPP)" ---
$ws1 = $wb.Worksheets.Item(1)

# B1: plot_HPP exact name placeholder
$ws1.Range("B1").Value = "[name 1]"

# Clear out all the single-value input cells (B2:B16) but keep their formatting
$ws1.Range("B2").ClearContents()
$ws1.Range("B3").ClearContents()
$ws1.Range("B4").ClearContents()
$ws1.Range("B5").ClearContents()
$ws1.Range("B6").ClearContents()
$ws1.Range("B7").ClearContents()
$ws1.Range("B8").ClearContents()
$ws1.Range("B9").ClearContents()
$ws1.Range("B10").ClearContents()
$ws1.Range("B11").ClearContents()
$ws1.Range("B12").ClearContents()
$ws1.Range("B13").ClearContents()
$ws1.Range("B14").ClearContents()
$ws1.Range("B15").ClearContents()
$ws1.Range("B16").ClearContents()

# Update selection to C3
$ws1.Range("C3").Select()

# --- Sheet 2: "Plot release rules (single HPP)" ---
$ws2 = $wb.Worksheets.Item(2)

# Remove the example values entirely (formatting cleared too, matching shrunk used range)
$ws2.Range("C1").Clear()
$ws2.Range("C2").Clear()
$ws2.Range("D2").Clear()

# Update selection to B11
$ws2.Range("B11").Select()

# --- Sheet 3: "Plot power output (multi HPP)" ---
$ws3 = $wb.Worksheets.Item(3)

# Replace example HPP names with generic placeholders
$ws3.Range("B2").Value = "[name 1]"
$ws3.Range("C2").Value = "[name 2]"
$ws3.Range("D2").Value = "[name 3]"

# Clear out the single-value input cells but keep their formatting
$ws3.Range("B4").ClearContents()
$ws3.Range("B5").ClearContents()
$ws3.Range("B6").ClearContents()
$ws3.Range("B7").ClearContents()
$ws3.Range("B8").ClearContents()
$ws3.Range("B9").ClearContents()
$ws3.Range("B10").ClearContents()

# Re-fit row heights now that wrapped text content has changed/shrunk
$ws3.Rows.Item(2).EntireRow.AutoFit()
$ws3.Rows.Item(9).EntireRow.AutoFit()

# Update selection to C8
$ws3.Range("C8").Select()

# Re-select sheet 1 as the active tab/cell when done
$ws1.Range("C3").Select()
